# Weekly update of "Espárragos" (asparagus) price data for Vega Monumental
# Concepción. Each data row's values (Fecha, Variedad, Volumen, Precio
# mínimo/máximo/promedio, Unidad de comercialización, Origen, Precio $/Kg)
# are rotated to new rows as the weekly data series advances. Row 8
# (2021-12-15) is untouched; rows 2,3,4,5,6,7,9,10,11 receive the values
# that, in the prior version of the sheet, lived in rows 5,10,2,3,9,6,11,4,7
# respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2  (was row 5: 2021-10-20)
$ws.Range("D2").Value = 44489
$ws.Range("K2").Value = 1400
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = 1450
$ws.Range("P2").Value = 1450

# Row 3  (was row 10: 2021-11-19)
$ws.Range("D3").Value = 44519
$ws.Range("J3").Value = 250
$ws.Range("K3").Value = 1200
$ws.Range("L3").Value = 1300
$ws.Range("M3").Value = 1240
$ws.Range("P3").Value = 1240

# Row 4  (was row 2: 2021-11-11)
$ws.Range("D4").Value = 44511
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 1300
$ws.Range("L4").Value = 1400
$ws.Range("M4").Value = 1350
$ws.Range("N4").Value = "$/kilo"
$ws.Range("P4").Value = 1350

# Row 5  (was row 3: 2021-11-10)
$ws.Range("D5").Value = 44510
$ws.Range("K5").Value = 1300
$ws.Range("L5").Value = 1400
$ws.Range("M5").Value = 1350
$ws.Range("P5").Value = 1350

# Row 6  (was row 9: 2021-09-29)
$ws.Range("D6").Value = 44468
$ws.Range("H6").Value = "Verde"
$ws.Range("J6").Value = 500
$ws.Range("K6").Value = 1800
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1920
$ws.Range("O6").Value = "Provincia de Linares"
$ws.Range("P6").Value = 1920

# Row 7  (was row 6: 2021-11-24)
$ws.Range("D7").Value = 44524
$ws.Range("J7").Value = 200
$ws.Range("O7").Value = "Provincia de Talca"

# Row 9  (was row 11: 2021-10-08)
$ws.Range("D9").Value = 44477
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("K9").Value = 1400
$ws.Range("L9").Value = 1500
$ws.Range("M9").Value = 1460
$ws.Range("P9").Value = 1460

# Row 10  (was row 4: 2021-10-27)
$ws.Range("D10").Value = 44496
$ws.Range("J10").Value = 550
$ws.Range("K10").Value = 1500
$ws.Range("L10").Value = 2000
$ws.Range("M10").Value = 1773
$ws.Range("N10").Value = "$/paquete"
$ws.Range("P10").Value = 1773

# Row 11  (was row 7: 2021-11-26)
$ws.Range("D11").Value = 44526
$ws.Range("J11").Value = 100
$ws.Range("K11").Value = 1500
$ws.Range("L11").Value = 1600
$ws.Range("M11").Value = 1550
$ws.Range("P11").Value = 1550
